$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleD = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.305.80"
$ws.Range("D2").Style = $styleD
$ws.Range("E2").Value = "  -0.01%  "
$styleD = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.568.94"
$ws.Range("D3").Style = $styleD
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.03%  "
$styleD = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.26"
$ws.Range("D5").Style = $styleD
$ws.Range("E5").Value = "  +3.70%  "
$styleD = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.27"
$ws.Range("D6").Style = $styleD
$ws.Range("E6").Value = "  -0.07%  "
$styleD = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.563.39"
$ws.Range("D7").Style = $styleD
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +7.50%  "
$styleD = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.645"
$ws.Range("D11").Style = $styleD
$ws.Range("E11").Value = "  -0.27%  "
$styleD = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.96"
$ws.Range("D12").Style = $styleD
$ws.Range("E12").Value = "  -1.08%  "
$styleD = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000310"
$ws.Range("D13").Style = $styleD
$ws.Range("E13").Value = "  +0.85%  "
$styleD = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.49"
$ws.Range("D14").Style = $styleD
$ws.Range("E14").Value = "  -0.14%  "
$styleD = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.133.52"
$ws.Range("D15").Style = $styleD
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$styleD = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.605.56"
$ws.Range("D16").Style = $styleD
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$styleD = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.383.78"
$ws.Range("D17").Style = $styleD
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$styleD = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.72"
$ws.Range("D18").Style = $styleD
$ws.Range("E18").Value = "  +1.83%  "
$styleD = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.00"
$ws.Range("D19").Style = $styleD
$ws.Range("E19").Value = "  -2.44%  "
$styleD = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "580.44"
$ws.Range("D20").Style = $styleD
$ws.Range("E20").Value = "  +6.67%  "
$ws.Range("E21").Value = "  +0.45%  "
$styleD = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("D22").Style = $styleD
$ws.Range("E22").Value = "  -2.30%  "
$styleD = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.36"
$ws.Range("D23").Style = $styleD
$ws.Range("E23").Value = "  -3.49%  "
$styleD = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.74"
$ws.Range("D24").Style = $styleD
$ws.Range("E24").Value = "  +1.04%  "
$styleD = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.85"
$ws.Range("D25").Style = $styleD
$ws.Range("E25").Value = "  -0.79%  "
$styleD = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.84"
$ws.Range("D26").Style = $styleD
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  -2.58%  "
$styleD = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("D29").Style = $styleD
$ws.Range("E29").Value = "  +2.87%  "
$styleD = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.31"
$ws.Range("D30").Style = $styleD
$ws.Range("E30").Value = "  +0.10%  "
$styleD = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("D31").Style = $styleD
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("E32").Value = "  -2.17%  "
$styleD = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("D33").Style = $styleD
$ws.Range("E33").Value = "  +0.02%  "
$styleD = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.60"
$ws.Range("D34").Style = $styleD
$ws.Range("E34").Value = "  -2.53%  "
$styleD = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.70"
$ws.Range("D35").Style = $styleD
$ws.Range("E35").Value = "  +19.85%  "
$styleD = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.25"
$ws.Range("D36").Style = $styleD
$ws.Range("E36").Value = "  +0.90%  "
$styleD = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "532.47"
$ws.Range("D37").Style = $styleD
$ws.Range("E37").Value = "  -4.21%  "
$styleD = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("D38").Style = $styleD
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("E39").Value = "  +0.09%  "
$styleD = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.23"
$ws.Range("D40").Style = $styleD
$ws.Range("E40").Value = "  -3.56%  "
$styleD = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0786"
$ws.Range("D41").Style = $styleD
$ws.Range("E41").Value = "  +2.37%  "
$styleD = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.528.20"
$ws.Range("D42").Style = $styleD
$ws.Range("E42").Value = "  +4.84%  "
$styleD = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.52"
$ws.Range("D43").Style = $styleD
$ws.Range("E43").Value = "  +3.54%  "
$styleD = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.136"
$ws.Range("D44").Style = $styleD
$ws.Range("E44").Value = "  +1.05%  "
$styleD = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0460"
$ws.Range("D45").Style = $styleD
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("E46").Value = "  -1.43%  "
$styleD = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.44"
$ws.Range("D47").Style = $styleD
$ws.Range("E47").Value = "  -2.95%  "
$styleD = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.140"
$ws.Range("D48").Style = $styleD
$ws.Range("E48").Value = "  +2.96%  "
$styleD = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.22"
$ws.Range("D49").Style = $styleD
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +0.22%  "
$styleD = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.45"
$ws.Range("D51").Style = $styleD
$ws.Range("E51").Value = "  -0.71%  "
